# Fix two typos in the "Maureen Simon Foods" row (row 22) of the
# supplier_markers worksheet:
#   - Ingredients (column D): curly "quoted" name -> (parenthesized) name
#   - Description (column E): stray trailing "w" typo removed

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("supplier_markers")

$ws.Range("E22").Value = "Maureen’s is a local food manufacturer in Vancouver, B.C. Canada. We create unique plant-based, local-ethnic, original food products. At Maureen’s we are based on cooking from ‘long ago’. Caribbean cooking traditionally uses a lot of veggie-focused dishes and the flavours we lean on do their best to bring out the best in veg."

$ws.Range("D22").Value = 'Tofu Scramble Roll - (Rolliis)'
